$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12. This shifts the existing rows 12..145
# down to 13..146 (matching the observed diff, where every data row from
# 13 to 145 takes on the values previously held by the row above it),
# and leaves a blank row 12 ready for the new weekly entry.
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with the new week's data: same template values
# as the rest of the sheet, but with its own date and volume.
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C12").Value = "Los Lagos"
$ws.Range("D12").Value = 44530
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 100112039
$ws.Range("G12").Value = "Ciboulette"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 320
$ws.Range("K12").Value = 2500
$ws.Range("L12").Value = 2500
$ws.Range("M12").Value = 2500
$ws.Range("N12").Value = "$/docena de atados"
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 833
$ws.Range("Q12").Value = 3
$ws.Range("R12").Value = "Hortaliza"

# Ensure row 12's date cell keeps the same date number format as the rest
# of column D.
$ws.Range("D12").NumberFormat = $ws.Range("D13").NumberFormat
